$d = $word.ActiveDocument

# Locate the "Delivery description: ${delivery_description}" paragraph and
# remember its 1-based index within Paragraphs so we can re-fetch its
# successor after inserting a new paragraph (inserted ranges/paragraphs can
# go stale once the document mutates).
$target = $null
$targetIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Delivery description*") {
        $target = $p
        $targetIdx = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Delivery description' paragraph"
}

# Add spacing-after = 0 to that paragraph (matches the diff's added
# <w:pPr><w:spacing w:after="0"/></w:pPr>)
$target.SpaceAfter = 0

# Insert a brand-new paragraph right after it.
$target.Range.InsertParagraphAfter()

# Re-fetch the newly created (now populated by position) paragraph and fill
# it in with the payment-method text, also with spacing-after = 0.
$newPara = $d.Paragraphs($targetIdx + 1)
$newPara.Range.Text = "Payment method: `${payment_method}"
$newPara.SpaceAfter = 0
